# Scheduled market-data refresh for the Exodus_Profits workbook.
# Updates currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ) and LeveProfit(NQ/HQ)
# columns (H, I, J, K, L, M, N) on each leve-item sheet with freshly
# pulled prices, recomputed from the same Leve Gil baseline.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 249.16667
$ws.Range("I5").Value = 249.16667
$ws.Range("K5").Value = 249.16667
$ws.Range("M5").Value = -134.16667
$ws.Range("H7").Value = 4000
$ws.Range("J7").Value = 4000
$ws.Range("L7").Value = 4000
$ws.Range("N7").Value = -4224
$ws.Range("H9").Value = 451.0909
$ws.Range("I9").Value = 106.666664
$ws.Range("J9").Value = 2001
$ws.Range("K9").Value = 106.666664
$ws.Range("L9").Value = 2001
$ws.Range("M9").Value = 62.333336
$ws.Range("N9").Value = -2339
$ws.Range("H10").Value = 10040.8
$ws.Range("J10").Value = 15666.667
$ws.Range("L10").Value = 15666.667
$ws.Range("N10").Value = -16252.667
$ws.Range("H14").Value = 4000
$ws.Range("J14").Value = 4000
$ws.Range("L14").Value = 4000
$ws.Range("N14").Value = -4382
$ws.Range("H43").Value = 2636.5
$ws.Range("I43").Value = 2366.3333
$ws.Range("J43").Value = 2798.6
$ws.Range("K43").Value = 2366.3333
$ws.Range("L43").Value = 2798.6
$ws.Range("M43").Value = -2297.3333
$ws.Range("N43").Value = -2936.6
$ws.Range("H64").Value = 7953.1816
$ws.Range("I64").Value = 7948.8
$ws.Range("J64").Value = 7997
$ws.Range("K64").Value = 7948.8
$ws.Range("L64").Value = 7997
$ws.Range("M64").Value = -7700.8
$ws.Range("N64").Value = -8493
$ws.Range("H67").Value = 7953.1816
$ws.Range("I67").Value = 7948.8
$ws.Range("J67").Value = 7997
$ws.Range("K67").Value = 7948.8
$ws.Range("L67").Value = 7997
$ws.Range("M67").Value = -7090.8
$ws.Range("N67").Value = -9713
$ws.Range("H86").Value = 2726.1292
$ws.Range("I86").Value = 2317.7896
$ws.Range("K86").Value = 2317.7896
$ws.Range("M86").Value = -1194.7896
$ws.Range("H89").Value = 2726.1292
$ws.Range("I89").Value = 2317.7896
$ws.Range("K89").Value = 11588.948
$ws.Range("M89").Value = -5972.948
$ws.Range("H111").Value = 669.6
$ws.Range("I111").Value = 407.66666
$ws.Range("K111").Value = 1222.99998
$ws.Range("M111").Value = 1844.00002
$ws.Range("H113").Value = 4116.8237
$ws.Range("I113").Value = 3890.75
$ws.Range("J113").Value = 4659.4
$ws.Range("K113").Value = 3890.75
$ws.Range("L113").Value = 4659.4
$ws.Range("M113").Value = -636.75
$ws.Range("N113").Value = -11167.4
$ws.Range("H125").Value = 5425.231
$ws.Range("J125").Value = 6003.778
$ws.Range("L125").Value = 54034.002
$ws.Range("N125").Value = -58954.002
$ws.Range("H132").Value = 1152.4584
$ws.Range("I132").Value = 1182.6364
$ws.Range("J132").Value = 820.5
$ws.Range("K132").Value = 3547.9092
$ws.Range("L132").Value = 2461.5
$ws.Range("M132").Value = -1017.9092
$ws.Range("N132").Value = -7521.5
$ws.Range("H137").Value = 339206.97
$ws.Range("I137").Value = 1390.4
$ws.Range("J137").Value = 1817154.5
$ws.Range("K137").Value = 4171.200000000001
$ws.Range("L137").Value = 5451463.5
$ws.Range("M137").Value = -1621.200000000001
$ws.Range("N137").Value = -5456563.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 42747.4
$ws.Range("I74").Value = 57303.945
$ws.Range("J74").Value = 5316.2856
$ws.Range("K74").Value = 57303.945
$ws.Range("L74").Value = 5316.2856
$ws.Range("M74").Value = -56429.945
$ws.Range("N74").Value = -7064.2856
$ws.Range("H77").Value = 42747.4
$ws.Range("I77").Value = 57303.945
$ws.Range("J77").Value = 5316.2856
$ws.Range("K77").Value = 286519.725
$ws.Range("L77").Value = 26581.428
$ws.Range("M77").Value = -282151.725
$ws.Range("N77").Value = -35317.428
$ws.Range("H132").Value = 1807.7878
$ws.Range("I132").Value = 1705.2667
$ws.Range("K132").Value = 5115.800099999999
$ws.Range("M132").Value = -2585.800099999999
$ws.Range("H138").Value = 122633.336
$ws.Range("J138").Value = 122633.336
$ws.Range("L138").Value = 122633.336
$ws.Range("N138").Value = -132913.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2349146
$ws.Range("I99").Value = 92511.63
$ws.Range("J99").Value = 5895285.5
$ws.Range("K99").Value = 92511.63
$ws.Range("L99").Value = 5895285.5
$ws.Range("M99").Value = -91013.63
$ws.Range("N99").Value = -5898281.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2927.0833
$ws.Range("I31").Value = 1825.6
$ws.Range("K31").Value = 1825.6
$ws.Range("M31").Value = -1530.6
$ws.Range("H34").Value = 2927.0833
$ws.Range("I34").Value = 1825.6
$ws.Range("K34").Value = 1825.6
$ws.Range("M34").Value = -1623.6
$ws.Range("H86").Value = 7147952
$ws.Range("I86").Value = 11908754
$ws.Range("J86").Value = 6748.5
$ws.Range("K86").Value = 11908754
$ws.Range("L86").Value = 6748.5
$ws.Range("M86").Value = -11907631
$ws.Range("N86").Value = -8994.5
$ws.Range("H89").Value = 7147952
$ws.Range("I89").Value = 11908754
$ws.Range("J89").Value = 6748.5
$ws.Range("K89").Value = 59543770
$ws.Range("L89").Value = 33742.5
$ws.Range("M89").Value = -59538154
$ws.Range("N89").Value = -44974.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 2678.8572
$ws.Range("I36").Value = 2350.4
$ws.Range("J36").Value = 3500
$ws.Range("K36").Value = 7051.200000000001
$ws.Range("L36").Value = 10500
$ws.Range("M36").Value = -6882.200000000001
$ws.Range("N36").Value = -10838

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3744.5715
$ws.Range("I126").Value = 2375
$ws.Range("J126").Value = 4292.4
$ws.Range("K126").Value = 7125
$ws.Range("L126").Value = 12877.2
$ws.Range("M126").Value = -4655
$ws.Range("N126").Value = -17817.2
$ws.Range("H132").Value = 4025.8276
$ws.Range("I132").Value = 3282.652
$ws.Range("J132").Value = 6874.6665
$ws.Range("K132").Value = 9847.956
$ws.Range("L132").Value = 20623.9995
$ws.Range("M132").Value = -7317.956
$ws.Range("N132").Value = -25683.9995
$ws.Range("H141").Value = 44047.2
$ws.Range("J141").Value = 37059
$ws.Range("L141").Value = 37059
$ws.Range("N141").Value = -47419

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1680.6471
$ws.Range("I93").Value = 1399
$ws.Range("K93").Value = 1399
$ws.Range("M93").Value = -151
$ws.Range("H132").Value = 2813.476
$ws.Range("I132").Value = 1940.3334
$ws.Range("J132").Value = 3977.6667
$ws.Range("K132").Value = 5821.0002
$ws.Range("L132").Value = 11933.0001
$ws.Range("M132").Value = -3291.0002
$ws.Range("N132").Value = -16993.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3124.5386
$ws.Range("I126").Value = 1952.375
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 5857.125
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -3387.125
$ws.Range("N126").Value = -19940
$ws.Range("H140").Value = 67740
$ws.Range("J140").Value = 67740
$ws.Range("L140").Value = 67740
$ws.Range("N140").Value = -78100
$ws.Range("H141").Value = 99949.5
$ws.Range("J141").Value = 99949.5
$ws.Range("L141").Value = 99949.5
$ws.Range("N141").Value = -110309.5

Write-Output "Applied 183 cell updates"
